$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.145.88'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '3.717.44'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("D4").Value = '2.45'
$ws.Range("E4").Value = '  +29.44%  '
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '231.20'
$ws.Range("E6").Value = '  -2.15%  '
$ws.Range("D7").Value = '658.69'
$ws.Range("E7").Value = '  +1.02%  '
$ws.Range("D8").Value = '0.453'
$ws.Range("E8").Value = '  +6.69%  '
$ws.Range("D9").Value = '1.16'
$ws.Range("E9").Value = '  +10.55%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '3.713.73'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '48.50'
$ws.Range("E12").Value = '  +9.26%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.211'
$ws.Range("E13").Value = '  +3.76%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = '0.0000309'
$ws.Range("E14").Value = '  +5.12%  '
$ws.Range("D15").Value = '6.76'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '4.410.80'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '97.064.31'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = '8.95'
$ws.Range("E18").Value = '  +15.51%  '
$ws.Range("D19").Value = '3.736.72'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("D20").Value = '19.39'
$ws.Range("E20").Value = '  +3.36%  '
$ws.Range("D21").Value = '13.16'
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").Value = '0.556'
$ws.Range("E22").Value = '  +10.32%  '
$ws.Range("D23").Value = '543.96'
$ws.Range("E23").Value = '  +4.88%  '
$ws.Range("D24").Value = '3.37'
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").Value = '122.53'
$ws.Range("E25").Value = '  +20.99%  '
$ws.Range("B26").Value = 'Hedera'
$ws.Range("C26").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D26").Value = '0.239'
$ws.Range("E26").Value = '  +42.06%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0000215'
$ws.Range("E27").Value = '  +5.27%  '
$ws.Range("D28").Value = '6.90'
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.918.28'
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").Value = '13.13'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '13.65'
$ws.Range("E31").Value = '  +13.09%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '3.06'
$ws.Range("E32").Value = '  +1.61%  '
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").Value = '0.190'
$ws.Range("E34").Value = '  +3.50%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '33.80'
$ws.Range("E35").Value = '  +4.85%  '
$ws.Range("D36").Value = '0.630'
$ws.Range("E36").Value = '  +7.36%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '1.84'
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("B38").Value = 'Binance-PegBSC-USD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '628.93'
$ws.Range("E39").Value = '  -3.50%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = '8.51'
$ws.Range("E41").Value = '  -2.78%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '7.23'
$ws.Range("E42").Value = '  +6.10%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = '0.166'
$ws.Range("E43").Value = '  +4.50%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0511'
$ws.Range("E44").Value = '  +14.38%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.492'
$ws.Range("E45").Value = '  +14.55%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = '2.04'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '40.70'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '0.983'
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '2.37'
$ws.Range("E49").Value = '  +4.51%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '9.12'
$ws.Range("E50").Value = '  +8.23%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = '23.56'
$ws.Range("E51").Value = '  -0.01%  '
